$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 28, shifting existing rows 28..98 down to 29..99.
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with its data. Columns A, B, C, E, F, G,
# H, I, N, Q, R carry the same values as the row that used to be at 28
# (now row 29), so just re-enter them; D, J, K, L, M, O, P are the new values.
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44690
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112030
$ws.Range("G28").Value = "Poroto granado"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 480
$ws.Range("K28").Value = 24500
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 24750
$ws.Range("N28").Value = "$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia del Elquí"
$ws.Range("P28").Value = 990
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
